# Update the "dSF" (column F) values for specific rows as part of the
# "repull data, push all data, mean calculation" update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value  = -12
$ws.Range("F8").Value  = 3
$ws.Range("F10").Value = -1
$ws.Range("F13").Value = 9
$ws.Range("F14").Value = -4
$ws.Range("F19").Value = 2
$ws.Range("F20").Value = -1
